$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 429; this shifts the existing rows 429:467
# down to 430:468 (carrying their formatting/values along), matching the
# dimension change from A1:R467 to A1:R468.
$ws.Rows.Item(429).Insert()

# Populate the newly inserted row 429 with the new record.
$ws.Cells.Item(429, 1).Value = 6
$ws.Cells.Item(429, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(429, 3).Value = "Metropolitana"
$ws.Cells.Item(429, 4).Value = 44769
$ws.Cells.Item(429, 5).Value = 13
$ws.Cells.Item(429, 6).Value = 100112043
$ws.Cells.Item(429, 7).Value = "Pepino ensalada"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 260
$ws.Cells.Item(429, 11).Value = 17000
$ws.Cells.Item(429, 12).Value = 18000
$ws.Cells.Item(429, 13).Value = 17615
$ws.Cells.Item(429, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(429, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(429, 16).Value = 294
$ws.Cells.Item(429, 17).Value = 60
$ws.Cells.Item(429, 18).Value = "Hortaliza"
